$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(25).Copy()
$ws.Rows(26).Insert(-4121)
Write-Host "inserted"
